$d = $word.ActiveDocument

function Insert-Xml($range, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1) HOME paragraph: split "aplicativo web" out with grammar-check proofErr markers ---
$find1 = $d.Content
$found1 = $find1.Find.Execute("aplicativo web oferece", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $find1.Paragraphs(1)
$rng1 = $d.Range($para1.Range.Start, $para1.Range.End - 1)
$body1 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">HOME – </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">resumo do que o </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>aplicativo web</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> oferece </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">que são </w:t></w:r>' +
    '<w:r><w:t>as avaliações</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> das</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> instituições, </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">cursos </w:t></w:r>' +
    '<w:r><w:t>e matérias para melhorar o ensino no país, pela visão de alunos já cursando e assim ajudando a futuros acadêmicos a escolher de acordo com a sua própria possibilidade.</w:t></w:r>' +
    '</w:p>'
Insert-Xml $rng1 $body1

# --- 2) "Sugestão 2" paragraph: remove underline pPr, rewrite text with spellcheck markers
#        and a bookmark, then add a trailing tab-only paragraph after it ---
$find2 = $d.Content
$found2 = $find2.Find.Execute("Sugestão 2:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $find2.Paragraphs(1)
$rng2 = $para2.Range
$body2 =
    '<w:p>' +
        '<w:r><w:tab/><w:t xml:space="preserve">Sugestão 2: </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">Remover telas de cadastro </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Usuario</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> e </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>visualizarUsuario</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> na página de aluno</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t xml:space="preserve">, colocar aviso de </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>email</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> ou senha incorretos na tela de login.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:r><w:tab/></w:r></w:p>'
Insert-Xml $rng2 $body2
